# This workbook's data rows (2-73) are being re-sorted: the full contents of
# each row (all columns A:AY) move as a unit to a new row position. The
# mapping below gives, for each destination (new) row number, which source
# (old) row number's data should end up there. It was derived by matching
# the unique "Id" values in column A between the old and new layouts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 2
$lastDataRow = 73

# Columns that hold plain text dates formatted like "2023-08-26"
# (Startdatum / Slutdatum). Writing these through Value2 normally makes
# Excel auto-detect them as real dates and reformat the cell, which the
# source file does not do (they are stored as literal text). We write
# those two columns separately, below, using a leading apostrophe to force
# text interpretation, then clear the resulting cell format so no stray
# number-format/style gets attached.
$dateTextCols = @(25, 27)  # Y, AA

# new row -> old row
$mapping = @{
    2=13; 3=32; 4=27; 5=54; 6=59; 7=69; 8=52; 9=53; 10=60; 11=18;
    12=63; 13=24; 14=62; 15=42; 16=7; 17=39; 18=40; 19=38; 20=56;
    21=44; 22=57; 23=22; 24=46; 25=51; 26=23; 27=34; 28=26; 29=12;
    30=14; 31=66; 32=17; 33=72; 34=67; 35=55; 36=5; 37=36; 38=43;
    39=19; 40=73; 41=29; 42=9; 43=3; 44=21; 45=71; 46=41; 47=6;
    48=47; 49=28; 50=25; 51=20; 52=2; 53=4; 54=10; 55=58; 56=50;
    57=11; 58=49; 59=30; 60=31; 61=65; 62=15; 63=33; 64=16; 65=37;
    66=45; 67=61; 68=70; 69=68; 70=48; 71=35; 72=64; 73=8
}

# 1) Snapshot every existing data row's values (all columns, A:AY) before
#    making any changes, so overwrites never clobber data still needed as
#    a source for a later destination row.
$snapshot = @{}
for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $rowRange = $ws.Range("A${r}:AY${r}")
    $snapshot[$r] = $rowRange.Value2
}

# 2) Write back each destination row using the snapshot of its source row,
#    for every column EXCEPT the plain-text date columns handled in step 3.
for ($newRow = $firstDataRow; $newRow -le $lastDataRow; $newRow++) {
    $oldRow = $mapping[$newRow]
    $srcVals = $snapshot[$oldRow]

    $destRange = $ws.Range("A${newRow}:AY${newRow}")
    $destVals = $destRange.Value2

    for ($col = 1; $col -le 51; $col++) {
        if ($dateTextCols -notcontains $col) {
            $destVals[1, $col] = $srcVals[1, $col]
        }
    }
    $destRange.Value2 = $destVals
}

# 3) Now handle the plain-text date columns (Y, AA) one cell at a time,
#    forcing text storage with a leading apostrophe and then clearing the
#    format Excel attaches for that, so the final cell is a plain text
#    value with the default (unset) style, matching the source layout.
foreach ($col in $dateTextCols) {
    for ($newRow = $firstDataRow; $newRow -le $lastDataRow; $newRow++) {
        $oldRow = $mapping[$newRow]
        $srcVals = $snapshot[$oldRow]
        $textVal = $srcVals[1, $col]
        $cell = $ws.Cells.Item($newRow, $col)
        if ($null -eq $textVal -or $textVal -eq "") {
            $cell.Value2 = $null
        } else {
            $cell.Value2 = "'" + $textVal
        }
    }
    $colLetter = $ws.Cells.Item(1, $col).Address($false, $false) -replace '[0-9]', ''
    $ws.Range("${colLetter}${firstDataRow}:${colLetter}${lastDataRow}").ClearFormats()
}
